$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 43 (shifts old rows 43.. down by 2)
$ws.Rows.Item(43).Insert()
$ws.Rows.Item(43).Insert()

# New row 43: Jengibre, "Primera", fecha 44571
$ws.Range("A43").Value = 9
$ws.Range("B43").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C43").Value = "Metropolitana"
$ws.Range("D43").Value = 44571
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = 100114007
$ws.Range("G43").Value = "Jengibre"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 610
$ws.Range("K43").Value = 12000
$ws.Range("L43").Value = 13000
$ws.Range("M43").Value = 12500
$ws.Range("N43").Value = "$/caja 13 kilos"
$ws.Range("O43").Value = "Perú"
$ws.Range("P43").Value = 962
$ws.Range("Q43").Value = 13
$ws.Range("R43").Value = "Hortaliza"

# New row 44: Jengibre, "Segunda", fecha 44571
$ws.Range("A44").Value = 9
$ws.Range("B44").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C44").Value = "Metropolitana"
$ws.Range("D44").Value = 44571
$ws.Range("E44").Value = 13
$ws.Range("F44").Value = 100114007
$ws.Range("G44").Value = "Jengibre"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Segunda"
$ws.Range("J44").Value = 106
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = 10000
$ws.Range("N44").Value = "$/caja 13 kilos"
$ws.Range("O44").Value = "Perú"
$ws.Range("P44").Value = 769
$ws.Range("Q44").Value = 13
$ws.Range("R44").Value = "Hortaliza"
